$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '33.834.16'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.779.08'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.85'
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.11'
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('E9').Value = '  +1.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0678'
$ws.Range('E10').Value = '  -5.26%  '
$ws.Range('E11').Value = '  +1.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.035.61'
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.26'
$ws.Range('E13').Value = '  +5.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.793.19'
$ws.Range('E14').Value = '  -1.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '33.866.93'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('E16').Value = '  -3.23%  '
$ws.Range('E17').Value = '  -1.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.62'
$ws.Range('E18').Value = '  -2.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '238.58'
$ws.Range('E19').Value = '  -3.31%  '
$ws.Range('E20').Value = '  -1.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('E22').Value = '  -2.91%  '
$ws.Range('E23').Value = '  -2.10%  '
$ws.Range('E24').Value = '  -2.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '160.27'
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.04'
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.09'
$ws.Range('E27').Value = '  -2.89%  '
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('E30').Value = '  +1.23%  '
$ws.Range('E31').Value = '  -2.63%  '
$ws.Range('E32').Value = '  -3.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.52'
$ws.Range('E33').Value = '  +0.54%  '
$ws.Range('E34').Value = '  -1.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.386.04'
$ws.Range('E35').Value = '  -2.03%  '
$ws.Range('E36').Value = '  -1.03%  '
$ws.Range('E37').Value = '  -1.81%  '
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('E39').Value = '  +2.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.25'
$ws.Range('E40').Value = '  +4.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '78.60'
$ws.Range('E42').Value = '  -3.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.58'
$ws.Range('E43').Value = '  +13.58%  '
$ws.Range('E44').Value = '  -3.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0₆0141'
$ws.Range('E45').Value = '  +13.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0506'
$ws.Range('E46').Value = '  +2.29%  '
$ws.Range('E47').Value = '  +3.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '107.01'
$ws.Range('E48').Value = '  -0.34%  '
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.937.00'
$ws.Range('E50').Value = '  -1.12%  '
